$d = $word.ActiveDocument

$replacements = @(
    @{Old = "Täydellisesti"; New = "Design: Täydellisesti"},
    @{Old = "Aikataulu"; New = "Design: Aikataulu"},
    @{Old = "Hyvät alihankkijat"; New = "Design: Hyvät alihankkijat"},
    @{Old = "Ei missään vaiheessa"; New = "Design: Ei missään vaiheessa"},
    @{Old = "Pelkkää voittoa"; New = "Design: Pelkkää voittoa"},
    @{Old = "Vähemmän virheitä kuvissa"; New = "Design: Vähemmän virheitä kuvissa"},
    @{Old = ":)"; New = "Design: :)"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
